# This script re-shuffles the per-row data (columns D, I, J, K, L, M, N, O, P, Q)
# for rows 2-16 on the active sheet according to a fixed row permutation,
# leaving columns A, B, C, E, F, G, H, R untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> source (old) row number
$rowMap = @{
    2  = 15
    3  = 16
    4  = 2
    5  = 10
    6  = 3
    7  = 9
    8  = 5
    9  = 6
    10 = 13
    11 = 14
    12 = 11
    13 = 4
    14 = 7
    15 = 12
    16 = 8
}

# Columns whose values move together with the row permutation.
$cols = @(4, 9, 10, 11, 12, 13, 14, 15, 16, 17)   # D, I, J, K, L, M, N, O, P, Q

# 1) Snapshot the "before" values for every relevant cell so that writes
#    to one row never clobber data that still needs to be read for another row.
$snapshot = @{}
foreach ($r in 2..16) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the permuted values back out.
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value = $src[$c]
    }
}
